$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Tube outer dimension "B" changed from 2" to 3" for every row (6-19).
#    "C" also grows from 2" to 3" but only for the rows that previously held
#    a 2" value (6-10); the 4" / 6" rows (11-19) are untouched.
# ---------------------------------------------------------------------------
$ws.Range("B6:B19").Value = 3
$ws.Range("C6:C10").Value = 3

# ---------------------------------------------------------------------------
# 2. New columns K:O - bending-stress / yield-strength check.
#    Headers (row 5) first - order matches the shared-string insertion order
#    of the source edit (M5, then K5, L5, N5, O5).
# ---------------------------------------------------------------------------
$ws.Range("M5").Value = "Yield"
$ws.Range("K5").Value = "M"
$ws.Range("L5").Value = "C"
$ws.Range("N5").Value = "Max KSI"
$ws.Range("O5").Value = "% strength"

# ---------------------------------------------------------------------------
# 3. Formulas for rows 6-19.
#    K = M (moment) = I/2 * H/2
#    L = C (half width) = C/2
#    M = Yield stress = (K*L)/F
#    N = Max KSI (material yield strength constant)
#    O = % strength = M/N   (shown as a percentage)
# ---------------------------------------------------------------------------
for ($r = 6; $r -le 19; $r++) {
    $ws.Range("K$r").Formula = "=I$r/2*H$r/2"
    $ws.Range("L$r").Formula = "=C$r/2"
    $ws.Range("M$r").Formula = "=(K$r*L$r)/F$r"
    $ws.Range("N$r").Value = 36000
    $ws.Range("O$r").Formula = "=(M$r/N$r)"
}

# Referencing H/I columns in the K/L/M formulas above makes the engine copy
# their custom number formats onto the new cells - put K:N back to the
# default "Normal" style (no explicit format) to match the source workbook.
$ws.Range("K6:N19").Style = "Normal"

# O column gets an explicit percentage number format.
$ws.Range("O6:O19").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# 4. Three more formatted (but empty) rows below the table, matching the
#    H/I/J number formats used throughout the table - as if the formatting
#    of row 19 had been extended downward.
# ---------------------------------------------------------------------------
$ws.Range("H19:J19").Copy() | Out-Null
$ws.Range("H20:J22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. Selection moves to G17.
# ---------------------------------------------------------------------------
$ws.Range("G17").Select() | Out-Null
